$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "305.75"
Set-TextValue "E2" "-0.56%"
Set-TextValue "G2" "21"

Set-TextValue "D3" "38.83"
Set-TextValue "E3" "6.87%"
Set-TextValue "G3" "21"

Set-TextValue "D4" "5.109"
Set-TextValue "E4" "0.88%"
Set-TextValue "G4" "21"

Set-TextValue "D5" "0.08080"
Set-TextValue "E5" "-0.43%"
Set-TextValue "G5" "21"

Set-TextValue "D6" "1.929"
Set-TextValue "E6" "-2.90%"
Set-TextValue "G6" "21"

Set-TextValue "D7" "4.192"
Set-TextValue "E7" "0.67%"
Set-TextValue "G7" "21"

Set-TextValue "E8" "1.51%"
Set-TextValue "G8" "21"

Set-TextValue "D9" "0.9308"
Set-TextValue "E9" "0.04%"
Set-TextValue "G9" "21"

Set-TextValue "D10" "0.1467"
Set-TextValue "E10" "1.17%"
Set-TextValue "G10" "21"

Set-TextValue "D11" "0.1926"
Set-TextValue "E11" "-0.26%"
Set-TextValue "G11" "21"

Set-TextValue "D12" "0.09067"
Set-TextValue "E12" "-0.98%"
Set-TextValue "G12" "21"

Set-TextValue "D13" "0.03527"
Set-TextValue "E13" "2.58%"
Set-TextValue "G13" "21"

Set-TextValue "D14" "0.09792"
Set-TextValue "E14" "-0.97%"
Set-TextValue "G14" "21"

Set-TextValue "D15" "0.001396"
Set-TextValue "E15" "-1.84%"
Set-TextValue "G15" "21"

Set-TextValue "D16" "0.005994"
Set-TextValue "E16" "-11.11%"
Set-TextValue "G16" "21"

Set-TextValue "D17" "3.777"
Set-TextValue "E17" "-1.46%"
Set-TextValue "G17" "21"

Set-TextValue "D18" "3.433"
Set-TextValue "E18" "0.05%"
Set-TextValue "G18" "21"

Set-TextValue "E19" "-0.42%"
Set-TextValue "G19" "21"

Set-TextValue "E20" "2.57%"
Set-TextValue "G20" "21"

Set-TextValue "E21" "-2.93%"
Set-TextValue "G21" "21"

Set-TextValue "D22" "0.2419"
Set-TextValue "E22" "3.24%"
Set-TextValue "G22" "21"

Set-TextValue "D23" "0.04383"
Set-TextValue "E23" "-0.24%"
Set-TextValue "G23" "21"

Set-TextValue "D24" "0.001239"
Set-TextValue "E24" "-0.08%"
Set-TextValue "G24" "21"

Set-TextValue "D25" "0.004271"
Set-TextValue "E25" "2.46%"
Set-TextValue "G25" "21"

Set-TextValue "D26" "0.0001304"
Set-TextValue "E26" "0.05%"
Set-TextValue "G26" "21"

Set-TextValue "G27" "21"

Set-TextValue "G28" "21"

Set-TextValue "G29" "21"

Set-TextValue "G30" "21"

Set-TextValue "G31" "21"

Set-TextValue "G32" "21"

Set-TextValue "G33" "21"

Set-TextValue "G34" "21"

Set-TextValue "G35" "21"

Set-TextValue "G36" "21"

Set-TextValue "G37" "21"

Set-TextValue "G38" "21"

Set-TextValue "D39" "0.02031"
Set-TextValue "E39" "-0.34%"
Set-TextValue "G39" "21"

Set-TextValue "D40" "0.05064"
Set-TextValue "E40" "-1.19%"
Set-TextValue "G40" "21"

Set-TextValue "D41" "0.007525"
Set-TextValue "E41" "0.52%"
Set-TextValue "G41" "21"

Set-TextValue "D42" "0.009752"
Set-TextValue "E42" "-4.02%"
Set-TextValue "G42" "21"

Set-TextValue "D43" "0.1343"
Set-TextValue "E43" "-2.10%"
Set-TextValue "G43" "21"

Set-TextValue "E44" "0.05%"
Set-TextValue "G44" "21"

Set-TextValue "D45" "0.009913"
Set-TextValue "E45" "0.50%"
Set-TextValue "G45" "21"

Set-TextValue "D46" "0.00006197"
Set-TextValue "E46" "-2.16%"
Set-TextValue "G46" "21"

Set-TextValue "E47" "-0.17%"
Set-TextValue "G47" "21"

Set-TextValue "G48" "21"

Set-TextValue "D49" "0.001805"
Set-TextValue "E49" "12.31%"
Set-TextValue "G49" "21"

Set-TextValue "D50" "0.00002106"
Set-TextValue "E50" "-0.17%"
Set-TextValue "G50" "21"

Set-TextValue "D51" "0.0002006"
Set-TextValue "E51" "-0.17%"
Set-TextValue "G51" "21"
